$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 11458
$ws.Range("F3").Value = 1970
$ws.Range("G3").Value = 65
$ws.Range("G4").Value = 0
$ws.Range("F5").Value = 858
$ws.Range("G5").Value = 70
$ws.Range("F6").Value = 2437
$ws.Range("G8").Value = 75
$ws.Range("F10").Value = 472
$ws.Range("F11").Value = 1381
$ws.Range("F13").Value = 137
$ws.Range("F15").Value = 1009
$ws.Range("F17").Value = 693
$ws.Range("F18").Value = 1152
$ws.Range("F21").Value = 20
$ws.Range("F24").Value = 326
$ws.Range("F26").Value = 271
$ws.Range("F27").Value = 482
$ws.Range("F28").Value = 512
$ws.Range("F29").Value = 696

$ws = $wb.Worksheets.Item("演出")
$ws.Range("G2").Value = 0
$ws.Range("F5").Value = 915
$ws.Range("F7").Value = 75
$ws.Range("F9").Value = 110
$ws.Range("F10").Value = 48
$ws.Range("F11").Value = 421
$ws.Range("F13").Value = 6

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 93

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 11458
$ws.Range("F3").Value = 1970
$ws.Range("G3").Value = 65
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("F6").Value = 858
$ws.Range("G6").Value = 70
$ws.Range("F7").Value = 2437
$ws.Range("G9").Value = 75
$ws.Range("F12").Value = 472
$ws.Range("F13").Value = 93
$ws.Range("F14").Value = 1381
$ws.Range("F17").Value = 137
$ws.Range("F18").Value = 915
$ws.Range("F20").Value = 1009
$ws.Range("F22").Value = 693
$ws.Range("F23").Value = 1152
$ws.Range("F26").Value = 20
$ws.Range("F29").Value = 326
$ws.Range("F31").Value = 75
$ws.Range("F33").Value = 271
$ws.Range("F35").Value = 110
$ws.Range("F36").Value = 110
$ws.Range("F37").Value = 482
$ws.Range("F38").Value = 512
$ws.Range("F39").Value = 696
$ws.Range("F41").Value = 48
$ws.Range("F43").Value = 421
$ws.Range("F46").Value = 6
